$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# Add new row 20 data to "Logs" sheet
$ws.Cells.Item(20, 1).Value = "Is er al nieuws?"
$ws.Cells.Item(20, 2).Value = "mailmind.test@zohomail.eu"
$ws.Cells.Item(20, 3).Value = "Testmail #10: Is er al nieuws?"
$ws.Cells.Item(20, 4).Value = "Overig"
$ws.Cells.Item(20, 5).Value = "Beste afzender,`nBedankt voor je e-mail. Om je vraag te beantwoorden hebben we wat meer context nodig. Zou je kunnen aangeven waarover je precies nieuws verwacht? Zo kunnen we je beter van dienst zijn.`nMet vriendelijke groet,`n[Naam bedrijf] E-mailassistent"
$ws.Cells.Item(20, 6).Value = "2025-08-03 14:53:24"
$ws.Cells.Item(20, 7).Value = "Ja"
$ws.Cells.Item(20, 8).Value = "Nee"
$ws.Cells.Item(20, 9).Value = "Ja"
$ws.Cells.Item(20, 10).Value = "Nee"
$ws.Rows.Item(20).AutoFit()

# Extend conditional formatting ranges to include the new row 20
$ws.Range("D2:D19").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D20"))
$ws.Range("G2:G19").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G20"))
$ws.Range("H2:H19").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H20"))
$ws.Range("I2:I19").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I20"))
$ws.Range("J2:J19").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J20"))

# Update Dashboard summary count for "Overig" category (row 3, col B)
$dash.Cells.Item(3, 2).Value = 6
